$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style of an existing header cell (H1) onto the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-15
$values = @{
    2  = @(7, 7)
    3  = @(8, 8)
    4  = @(8, 9)
    5  = @(6, 6)
    6  = @(6, 7)
    7  = @(6, 8)
    8  = @(9, 9)
    9  = @(8, 9)
    10 = @(6, 8)
    11 = @(9, 9)
    12 = @(8, 8)
    13 = @(8, 9)
    14 = @(7, 8)
    15 = @(1, 4)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
